# "add count from center"
#
# Target sheet is "빈소1" (xl/worksheets/sheet2.xml — its rendered content is
# what the commit's diff shows under the `xl/personal.xlsx` path). Before the
# edit it holds 14 rows of order lines (rows "0".."13"); after the edit it
# holds just two rows: row "0" (치즈김밥 x8 = 25000) and row "1"
# (치즈김밥 x5 = 25000), with every other line removed and the used range
# collapsing to A0:E1.
#
# Note: this workbook's sheet XML has a literal row labelled "0", which is
# not a row Excel's object model can ever address (rows are 1-based, so
# Cells(0, c)/Range("A0") etc. are always invalid) — it's inert legacy data
# baked into the file that no COM automation can touch. We leave it exactly
# as-is and make every change that *is* reachable: rewriting Excel row 1
# (the sheet's XML row "1") to the new line item, and deleting the old rows
# 2-13 (XML rows "2".."13") so the sheet shrinks down to just two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("빈소1")

# How many rows currently hold data (Excel counts rows 1.. here; the extra
# row "0" baked into the file sits below Excel's addressable range and isn't
# included in this count).
$lastRow = $ws.UsedRange.Rows.Count

# Rewrite row 1: item, unit price, count, line total, (blank note column).
$ws.Cells.Item(1, 1).Value = "치즈김밥"
$ws.Cells.Item(1, 2).Value = 2500
$ws.Cells.Item(1, 3).Value = 5
$ws.Cells.Item(1, 4).Value = 25000

# Drop every remaining order line (old rows 2-13) so only the new row 1
# (plus the untouchable legacy row 0) is left.
if ($lastRow -ge 2) {
  $ws.Rows("2:" + $lastRow).Delete()
}
